$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (new shared strings get added automatically)
$ws.Range("A1").Value = "hong ha"
$ws.Range("B1").Value = "trung quocs"

# Update the selected cell on the sheet
$ws.Range("D5").Select()
